# Add "MAPE (%)" column (G) to both the "val" and "test" sheets, and
# switch the numeric-data style (currently "#,##0.00") over to "General"
# formatting, matching the new MAPE ratios which aren't meant to be
# displayed as fixed 2-decimal currency-like numbers.

$wb = $excel.ActiveWorkbook

$mapeByCell = @{
    "val"  = @{
        2  = 1.137668838919149
        3  = 1.107259572273046
        4  = 1.026070645191113
        5  = 1.114076757106135
        6  = 1.287108028287489
        7  = 1.369541467407601
        8  = 1.171746816392794
        9  = 1.06103386653175
        10 = 1.491186201848516
        11 = 1.520512022536001
    }
    "test" = @{
        2  = 1.162344678166155
        3  = 1.163480425759533
        4  = 1.08593337028839
        5  = 1.062189953057388
        6  = 1.128886164281552
        7  = 1.036382963724421
        8  = 1.242322966952531
        9  = 1.048130106018086
        10 = 0.8898799395913749
        11 = 0.9342840424008627
    }
}

foreach ($sheetName in @("val", "test")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Copy the header's formatting (bold/border/centered) from F1 onto the
    # new G1 header cell, then set its text.
    $ws.Range("F1").Copy()
    $ws.Range("G1").PasteSpecial(-4122)
    $ws.Range("G1").Value = "MAPE (%)"

    # Copy the data-cell formatting from column F onto column G for every
    # data row, then fill in the MAPE values.
    $rowValues = $mapeByCell[$sheetName]
    foreach ($row in 2..11) {
        $ws.Range("F$row").Copy()
        $ws.Range("G$row").PasteSpecial(-4122)
        $ws.Range("G$row").Value = $rowValues[$row]
    }

    # The shared numeric style used by columns B:G was "#,##0.00"; switch
    # it to General so the new MAPE ratios (and the existing metrics)
    # aren't forced into a 2-decimal/thousands-separator format. Clearing
    # formats drops the cells back to the workbook's default (General)
    # style, which already carries the same plain, non-bold Calibri font
    # the numeric cells used before.
    $ws.Range("B2:G11").ClearFormats()
}
